# Update cryptocurrency price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column cells we are about to rewrite to Text format so
# values like "61.320.52" / "412.46" are stored verbatim instead of being
# auto-parsed into numbers by Excel.
$dCells = @("D2","D3","D5","D6","D7","D8","D10","D11","D12","D14","D15","D16","D17","D18","D21","D22","D23","D24","D25","D27","D28","D30","D32","D33","D35","D36","D38","D39","D40","D41","D42","D44","D45","D46","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '61.320.52'
$ws.Range("E2").Value = '  +7.71%  '

$ws.Range("D3").Value = '3.399.06'
$ws.Range("E3").Value = '  +4.78%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '412.46'
$ws.Range("E5").Value = '  +4.03%  '

$ws.Range("D6").Value = '121.78'
$ws.Range("E6").Value = '  +13.15%  '

$ws.Range("D7").Value = '3.395.63'
$ws.Range("E7").Value = '  +4.75%  '

$ws.Range("D8").Value = '0.575'
$ws.Range("E8").Value = '  -0.83%  '

$ws.Range("E9").Value = '  +0.04%  '

$ws.Range("D10").Value = '0.637'
$ws.Range("E10").Value = '  +3.09%  '

$ws.Range("D11").Value = '0.119'
$ws.Range("E11").Value = '  +25.17%  '

$ws.Range("D12").Value = '40.96'
$ws.Range("E12").Value = '  +4.76%  '

$ws.Range("E13").Value = '  -0.75%  '

$ws.Range("D14").Value = '3.942.14'
$ws.Range("E14").Value = '  +4.98%  '

$ws.Range("D15").Value = '8.36'
$ws.Range("E15").Value = '  +0.89%  '

$ws.Range("D16").Value = '19.48'
$ws.Range("E16").Value = '  +3.35%  '

$ws.Range("D17").Value = '3.395.59'
$ws.Range("E17").Value = '  +4.68%  '

$ws.Range("D18").Value = '61.311.15'
$ws.Range("E18").Value = '  +8.06%  '

$ws.Range("E19").Value = '  -0.97%  '

$ws.Range("E20").Value = '  -0.89%  '

$ws.Range("D21").Value = '0.0000122'
$ws.Range("E21").Value = '  +11.54%  '

$ws.Range("D22").Value = '3.32'
$ws.Range("E22").Value = '  -0.44%  '

$ws.Range("D23").Value = '12.80'
$ws.Range("E23").Value = '  -0.22%  '

$ws.Range("D24").Value = '297.62'
$ws.Range("E24").Value = '  +2.12%  '

$ws.Range("D25").Value = '76.10'
$ws.Range("E25").Value = '  +2.79%  '

$ws.Range("E26").Value = '  -1.05%  '

$ws.Range("D27").Value = '30.56'
$ws.Range("E27").Value = '  +9.17%  '

$ws.Range("D28").Value = '8.07'
$ws.Range("E28").Value = '  +11.63%  '

$ws.Range("E29").Value = '  -2.47%  '

$ws.Range("D30").Value = '7.60'
$ws.Range("E30").Value = '  -5.27%  '

$ws.Range("E31").Value = '  +1.75%  '

$ws.Range("D32").Value = '0.115'
$ws.Range("E32").Value = '  +4.82%  '

$ws.Range("D33").Value = '42.36'
$ws.Range("E33").Value = '  +3.80%  '

$ws.Range("E34").Value = '  +19.50%  '

$ws.Range("B35").Value = 'Dai'
$ws.Range("C35").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.07%  '

$ws.Range("B36").Value = 'Cosmos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D36").Value = '11.34'
$ws.Range("E36").Value = '  +1.56%  '

$ws.Range("E37").Value = '  -0.77%  '

$ws.Range("D38").Value = '52.39'
$ws.Range("E38").Value = '  +2.39%  '

$ws.Range("D39").Value = '3.51'
$ws.Range("E39").Value = '  +1.64%  '

$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  -0.08%  '

$ws.Range("D41").Value = '3.01'
$ws.Range("E41").Value = '  +1.47%  '

$ws.Range("D42").Value = '1.96'
$ws.Range("E42").Value = '  +5.09%  '

$ws.Range("E43").Value = '  +0.73%  '

$ws.Range("D44").Value = '133.29'
$ws.Range("E44").Value = '  -2.98%  '

$ws.Range("D45").Value = '17.15'
$ws.Range("E45").Value = '  +3.28%  '

$ws.Range("D46").Value = '3.91'
$ws.Range("E46").Value = '  +0.03%  '

$ws.Range("E47").Value = '  -0.16%  '

$ws.Range("E48").Value = '  -2.04%  '

$ws.Range("D49").Value = '21.69'
$ws.Range("E49").Value = '  -2.92%  '

$ws.Range("D50").Value = '2.200.35'
$ws.Range("E50").Value = '  +2.28%  '

$ws.Range("D51").Value = '3.744.44'
$ws.Range("E51").Value = '  +5.08%  '

# Restore the Normal style on those Price cells now that the text value is
# locked in, so no stray number-format style lingers on the cell.
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
